$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9574000239372253
$ws.Range("B1").Value = 1.35305655002594
$ws.Range("C1").Value = 2.443425416946411
$ws.Range("D1").Value = 4.143021106719971
$ws.Range("E1").Value = 1.859736323356628
